# daily auto push: 2026-02-24 19:16 UTC
#
# Two new readings for 2026/02/24 (23:00) and 2026/02/25 (01:00) were
# recorded. They belong right after the existing 2026/02/24 rows (the
# sheet is otherwise in chronological order up to that point, followed by
# a separate block that jumps ahead to 2026/12/29 onward), so insert two
# fresh rows at row 877 and push everything from the old row 877 down by
# two (917 -> 919, 918 -> 920, etc.), then populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 877, pushing existing data down.
$ws.Rows.Item(877).Insert()
$ws.Rows.Item(877).Insert()

# Make sure the date-like text in column A/B is stored as plain text,
# not auto-converted into a date serial number.
$ws.Range("A877:B878").NumberFormat = "@"

$ws.Range("A877").Value = "2026/02/24"
$ws.Range("B877").Value = "火"
$ws.Range("C877").Value = 23
$ws.Range("D877").Value = 201

$ws.Range("A878").Value = "2026/02/25"
$ws.Range("B878").Value = "水"
$ws.Range("C878").Value = 1
$ws.Range("D878").Value = 201

# Drop the formatting we applied above so the new rows match the plain
# (unstyled) look of the other data rows in the sheet.
$ws.Range("A877:D878").ClearFormats()
